$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write the standard 4-column header row ("line_n","prev_line","line",
# "next_line") onto a freshly added worksheet, bold + centered, matching the
# style already used by the other sheets in this workbook.
# ---------------------------------------------------------------------------
function Set-Header($ws) {
    $ws.Range("A1").Value = "line_n"
    $ws.Range("B1").Value = "prev_line"
    $ws.Range("C1").Value = "line"
    $ws.Range("D1").Value = "next_line"
    $hdr = $ws.Range("A1:D1")
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 1) Insert brand-new sheet "A05 vie sainte dieudonnee" right after
#    "A04 vie saint christofle" (currently sheet index 4).
# ---------------------------------------------------------------------------
$afterA04 = $wb.Worksheets.Item(4)
$a05 = $wb.Worksheets.Add($null, $afterA04)
$a05.Name = "A05 vie sainte dieudonnee"
Set-Header $a05
$a05.Range("A2").Value = 436
$a05.Range("B2").Value = "Ma suer, bien vous en croi quant me proumetés tant,"
$a05.Range("C2").Value = "Si parleray a li.” Lors s'en parti atant."
$a05.Range("D2").Value = "Viers l'evesque se traist, qui vint apriés diner"

# ---------------------------------------------------------------------------
# 2) Insert brand-new (blank) sheet "A10 poines enfer" right after
#    "A09 vie saint leu" (currently sheet index 9, after step 1's insert).
# ---------------------------------------------------------------------------
$afterA09 = $wb.Worksheets.Item(9)
$a10 = $wb.Worksheets.Add($null, $afterA09)
$a10.Name = "A10 poines enfer"
Set-Header $a10

# ---------------------------------------------------------------------------
# 3) "A16 guillaume angleterre" now sits at index 14. Insert a new data row
#    at row 2 (line_n 115) and correct the punctuation/wording throughout.
# ---------------------------------------------------------------------------
$a16 = $wb.Worksheets.Item(14)

$a16.Rows.Item(2).Insert()
$a16.Range("A2:D2").ClearFormats()

$a16.Range("A2").Value = 115
$a16.Range("B2").Value = "Ne vous ne vostre pere (dont Dieus asoille l'ame!),"
$a16.Range("C2").Value = "Que vous le fassiez rendre.” Cil qui fu sans disfame,"
$a16.Range("D2").Value = "Li dist: “Je le weil bien, foy que doy Nostre Dame!”"

$a16.Range("A3").Value = 169
$a16.Range("B3").Value = "Avras et grans meschief; mais Dieus, qui sur tous vaut,"
$a16.Range("C3").Value = "Le te rendra moult bien.” Lors la vois s'en parti;"
$a16.Range("D3").Value = "Et le roy d'Angleterre requist a Dieu merci,"

$a16.Range("A4").Value = 286
$a16.Range("B4").Value = "Que l'un de mes enfans mengier me covendra,"
$a16.Range("C4").Value = "Se je n'ai char ou pain.” Lors le roy souspira;"
$a16.Range("D4").Value = "Il a traite s'espee, sa chausce desferma"

$a16.Range("A5").Value = 295
$a16.Range("B5").Value = "Se je vous voi couper. Ja, se Dieu plest, mes dens"
$a16.Range("C5").Value = "N'usseront vostre char.” Le roi, qui fu dolens,"
$a16.Range("D5").Value = "Respondi: “Douce suer, ne puis trouver pourpens"

$a16.Range("A6").Value = 444
$a16.Range("B6").Value = "Forment prist a pleurer et dist: “Vierge Marie,"
$a16.Range("C6").Value = "Or ai je tout perdu.” Lors vint celle partie"
$a16.Range("D6").Value = "Ou la bource au marcheant fu pendant demouree,"

$a16.Range("A7").Value = 592
$a16.Range("B7").Value = "Doucement dist: “Seingneur, el non du dous Jhesus,"
$a16.Range("C7").Value = "Metez moi en vos nef.” L'un d'eus a respondu:"
$a16.Range("D7").Value = "“Biaus amis, dont viens tu? Moult as sousfert de paine;"

$a16.Range("A8").Value = 686
$a16.Range("B8").Value = "Je vous ai trop cousté: Dieu me doint vivre tant"
$a16.Range("C8").Value = "Que le vous puisse rendre.” Lors le preudon pleurant"
$a16.Range("D8").Value = "Li fist donner dix livres et un cheval courant;"

$a16.Range("A9").Value = 718
$a16.Range("B9").Value = "Le fortier dist: “Sire, bersant vont par le bois;"
$a16.Range("C9").Value = "Ceste beste ont tuee.” Le seingneur fu courtois,"
$a16.Range("D9").Value = "Il dist: “Je leur pardoins, pour Dieu, a ceste fois;"

$a16.Range("A10").Value = 797
$a16.Range("B10").Value = "Se le roy revenoit de quoy vous me parlés,"
$a16.Range("C10").Value = "Seroit il receus?” Son neveu dist errant:"
$a16.Range("D10").Value = "“Ouïl, se Dieus me gart, a sollanpnité grant;"

$a16.Range("A11").Value = 824
$a16.Range("B11").Value = "Cëans a tel jouel qui vault, je vous afi,"
$a16.Range("C11").Value = "Plus de tiex trente anniaus.” La dame respondi:"
$a16.Range("D11").Value = "“Je ne weil que l'annel.” Adont le roy li tant;"

# ---------------------------------------------------------------------------
# 4) "A17 robert deable" now sits at index 15. Correct punctuation/wording
#    in its three existing data rows (no row count change).
# ---------------------------------------------------------------------------
$a17 = $wb.Worksheets.Item(15)

$a17.Range("C2").Value = "Ou que ne m'apartiengne.” L'un d'euz dist sanz espasse:"
$a17.Range("D2").Value = "“Chier sire, je sai telle dont ja n'arez vergoigne;"

$a17.Range("B3").Value = "Qui de cuer le verroie morir devant mes yex,"
$a17.Range("C3").Value = "Trop m'a fait de courous.” Lors parla uns hons viex"
$a17.Range("D3").Value = "Et dist: “Sire, je lo que Robert soit mandez"

$a17.Range("B4").Value = "La duchesse dist: “Filz, certes la coupe est moie"
$a17.Range("C4").Value = "De vostre mauvaistie.” Lors li conta la voie"
$a17.Range("D4").Value = "Comment elle l'avoit au deable donné"
